# Auto-generated edit script: updates crypto price/volume table cells
# to match the "Updated cryptos list" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.269.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.63%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.576.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.98%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.35%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -0.52%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -1.90%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D8").Value = "'22.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.13%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.249"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.48%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.01%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.22%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.800.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.13%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.573.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.01%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -1.54%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.520"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.66%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'Litecoin"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'62.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.97%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'WrappedBTC"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'27.280.45"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.59%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'215.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.02%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.44%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  -1.04%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.30%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.34%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.57%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  +1.40%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'151.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.54%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  -5.61%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'14.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.05%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -1.34%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.35%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.85%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -2.30%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.11%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.404.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.66%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -1.73%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  +1.51%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("B36").Value = "'HuobiToken"
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").Value = "'2.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -2.39%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("B37").Value = "'TrustWalletToken"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.941"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.91%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.29%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.87%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -2.75%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +1.65%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.22%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +1.57%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.54%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'63.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.83%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.712.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.14%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'86.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.24%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.0₇0991"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.44%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0955"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.30%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.40%  "
$ws.Range("E51").Style = "Normal"
